# initialize trade testing via RSI, 1 share
#
# - watchlist: GOOG -> GOOGL
# - portfolio: STOCKS value initialized to 0 (1 share not bought yet, pending
#   the RSI signal)
# - (stocks sheet header cells are untouched - their shared-string indices
#   only shift because GOOG is replaced/reindexed in the shared string table,
#   the header text/order itself does not change)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# watchlist: GOOG -> GOOGL
# ---------------------------------------------------------------------------
$wsWatch = $wb.Worksheets.Item("watchlist")
$wsWatch.Range("A4").Value = "GOOGL"

# ---------------------------------------------------------------------------
# stocks: no content change, view only
# ---------------------------------------------------------------------------
$wsStocks = $wb.Worksheets.Item("stocks")

# ---------------------------------------------------------------------------
# portfolio: STOCKS starts at 0 value (before any RSI-triggered buys)
# ---------------------------------------------------------------------------
$wsPortfolio = $wb.Worksheets.Item("portfolio")
$wsPortfolio.Range("B3").Value = 0

# ---------------------------------------------------------------------------
# view state: zoom + selection per sheet, finishing on portfolio (active tab)
# ---------------------------------------------------------------------------
$wsWatch.Activate()
$wsWatch.Range("D35").Select()
$excel.ActiveWindow.Zoom = 175

$wsStocks.Activate()
$excel.ActiveWindow.Zoom = 130

$wsPortfolio.Activate()
$wsPortfolio.Range("C7").Select()
$excel.ActiveWindow.Zoom = 145
